# "Desenvolvido nova tela de adicao de turmas" - rework of the "add turma"
# sheet layout: widen the existing columns, add a new 6th column, fix the
# "PERIODO" header's accentuation/casing, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1 holds the shared string "PERIODO (Manhã/Tarde)" -> normalize to
# "PERIODO (MANHÃ/TARDE)". Re-assigning the same cell's Value updates the
# shared-string table in place instead of touching any other cell.
$ws.Range("E1").Value = "PERIODO (MANHÃ/TARDE)"

# New columns layout for the redesigned "add turma" screen: columns A-E are
# widened and a brand-new column F is introduced.
$ws.Columns.Item(1).ColumnWidth = 43.666666666666664
$ws.Columns.Item(2).ColumnWidth = 36.333333333333336
$ws.Columns.Item(3).ColumnWidth = 22.166666666666668
$ws.Columns.Item(4).ColumnWidth = 24.166666666666668
$ws.Columns.Item(5).ColumnWidth = 23
$ws.Columns.Item(6).ColumnWidth = 16

# The saved view now has the cursor sitting on E2 instead of A2.
$ws.Range("E2").Select()
